$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 3 for the "py" command, shifting existing rows down.
$ws.Rows(3).Insert()

# New row 3: py command
$ws.Range("A3").Value = "py"
$ws.Range("B3").Value = "Run python from a file"
$ws.Range("C3").Value = "py <filename>"
$ws.Range("D3").Value = "The shell object will be a global called 'shell' in the context of the Python script."

# Notes for the "source" command (row 2)
$ws.Range("D2").Value = "Commands are simply interpreted by the xcpshell."

# Update selection/view to match target state
$ws.Range("D2").Select()
